$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 173659
$ws.Range("C4").Value = 163879
$ws.Range("C5").Value = 9780
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 5.63
